$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated sample salinity for WS: column U (salinity) changes from 35 -> 32
# for every row whose salinity was 35 (rows 2-8 and 12-21; rows 9-11 used a
# different measured salinity and are left untouched).
$rows = @(2,3,4,5,6,7,8,12,13,14,15,16,17,18,19,20,21)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 21).Value = 32
}

# Reflect the author's on-screen scroll position / selection at save time.
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("U27:V27").Select()
